# Weekly update: a new daily price record for "Achicoria" (Vega Central
# Mapocho de Santiago) is inserted at the top of the data block (row 13,
# right after the header + first 11 existing records), pushing every
# existing record down by one row. The former last record (old row 73)
# becomes the new last record (row 74), growing the used range from
# A1:R73 to A1:R74.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 13 - this shifts rows 13:73 down to 14:74
# and extends the sheet's dimension automatically.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new record.
$ws.Range("A13").Value = 9
$ws.Range("B13").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C13").Value = "Metropolitana"
$ws.Range("D13").Value = 45114
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 100112010
$ws.Range("G13").Value = "Achicoria"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 90
$ws.Range("K13").Value = 8000
$ws.Range("L13").Value = 8000
$ws.Range("M13").Value = 8000
$ws.Range("N13").Value = "$/caja 16 unidades"
$ws.Range("O13").Value = "Provincia de Quillota"
$ws.Range("P13").Value = 500
$ws.Range("Q13").Value = 16
$ws.Range("R13").Value = "Hortaliza"
